$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: add a space before the ":" in the status labels ---
$ws.Range("B2").Value = '4 : pas de résultats postés ni publiés'
$ws.Range("B3").Value = '1 : résultats postés ou publiés dans les 12 mois'
$ws.Range("B4").Value = '4 : pas de résultats postés ni publiés'
$ws.Range("B5").Value = '4 : pas de résultats postés ni publiés'
$ws.Range("B6").Value = '4 : pas de résultats postés ni publiés'
$ws.Range("B7").Value = '2 : résultats postés ou publiés entre 12 et 36 mois'
$ws.Range("B8").Value = '4 : pas de résultats postés ni publiés'
$ws.Range("B9").Value = '4 : pas de résultats postés ni publiés'
$ws.Range("B10").Value = '2 : résultats postés ou publiés entre 12 et 36 mois'
$ws.Range("B11").Value = '3 : résultats postés ou publiés après les 36 mois'
$ws.Range("B12").Value = '4 : pas de résultats postés ni publiés'
$ws.Range("B13").Value = '4 : pas de résultats postés ni publiés'
$ws.Range("B14").Value = '4 : pas de résultats postés ni publiés'
$ws.Range("B15").Value = '4 : pas de résultats postés ni publiés'
$ws.Range("B16").Value = '4 : pas de résultats postés ni publiés'

# --- Column A: re-point the "statut" code to match the re-ordered row ---
# (copy from cells that already hold the right text so the value stays a
# text entry rather than being re-interpreted as a number by Excel)
$ws.Range("A10").Copy($ws.Range("A11"))
$ws.Range("A8").Copy($ws.Range("A7"))
$ws.Range("A8").Copy($ws.Range("A10"))
$ws.Range("A2").Copy($ws.Range("A8"))
$ws.Range("A2").Copy($ws.Range("A13"))

# --- Column C: NCTId re-ordering ---
$ws.Range("C7").Value = 'NCT04634318'
$ws.Range("C8").Value = 'NCT05237050'
$ws.Range("C9").Value = 'NCT05077605'
$ws.Range("C10").Value = 'NCT04028973'
$ws.Range("C11").Value = 'NCT05437991'
$ws.Range("C13").Value = 'NCT04878263'

# --- Column G: clinical_trial_title re-ordering ---
$ws.Range("G7").Value = 'Organization of Pulmonary Rehabilitation of Post-COVID-19 Patient With Sequelae. Assessment and Therapeutic Indication of Tele-rehabilitation Versus Conventional Rehabilitation'
$ws.Range("G8").Value = 'Evaluation of Sound Therapy in a Population of Women With Fibromyalgia Aged Between 30 and 60 Years'
$ws.Range("G9").Value = 'Electrical Impedance Tomography: Effect of Extubation on Functional Residual Capacity'
$ws.Range("G10").Value = 'Towards a Better Understanding of Neuromuscular Alterations and Fatigue in Chronic Obstructive Pulmonary Disease (COPD)'
$ws.Range("G11").Value = 'Ultrasonographic Morphology Assessment of Low-grade Carotid Stenosis'
$ws.Range("G13").Value = 'Evaluation and Support Care Process Within the Care Pathway of Heart Failure Patients'

# --- Column H: acronym re-ordering ---
$ws.Range("H7").Value = 'REHABCOVID'
$ws.Range("H8").Value = 'SONOMYAL'
$ws.Range("H9").Value = 'PULMOVISTA'
$ws.Range("H10").Value = 'NEUROTIGUE'
$ws.Range("H11").Value = 'QUAMUS'
$ws.Range("H13").Value = 'FIL-EAS'

# --- Column I: intervention_type swap ---
$ws.Range("I10").Value = 'OTHER'
$ws.Range("I11").Value = 'DIAGNOSTIC_TEST'
